$d = $word.ActiveDocument

# --- Merge four "name -> blank line -> paragraph" pairs into "name -> paragraph"
# by deleting the blank paragraph's paragraph mark (joins it with the next
# paragraph, which is exactly what the diff shows: the empty <w:p> plus the
# following paragraph's opening <w:p><w:pPr> collapse into one paragraph).
# Work from the bottom of the document upward so earlier paragraph indices
# stay valid as later ones are removed.
$d.Paragraphs.Item(33).Range.Delete()   # blank line before "Broadly, I am interested..."
$d.Paragraphs.Item(29).Range.Delete()   # blank line before "With limited access..."
$d.Paragraphs.Item(24).Range.Delete()   # blank line before "In my dissertation..."
$d.Paragraphs.Item(20).Range.Delete()   # blank line before "My research presents..."

# --- Append the two new paragraphs (Kristin's passage) right after the
# "...how to live meaningfully." paragraph, which is now paragraph 30.
$pBroadly = $d.Paragraphs.Item(30)
$pBroadly.Range.InsertParagraphAfter()

$pKristin = $d.Paragraphs.Item(31)
$pKristin.Range.InsertAfter("Kristin")
$kristinStart = $pKristin.Range
$kristinStart.SetRange($kristinStart.Start, $kristinStart.Start)
$kristinStart.InsertAfter([char]11)

$pKristin = $d.Paragraphs.Item(31)
$pKristin.Range.InsertParagraphAfter()

$pMusic = $d.Paragraphs.Item(32)
$musicText = "Generally, music composers are concerned with a musical experience, rather than the revelation of compositional materials. But when the source of the data is made explicit, it raises the question of whether some aspect of the source phenomenon can be understood by listening to the piece. When the primary intention of the composer shifts to the revelation of the source, the work crosses into the realm of sonification. With this crossing over comes a question of whether the listener can also understand the composer’s intention to produce more than an experience of the music itself."
$pMusic.Range.InsertAfter($musicText)

# --- Narrow the left/right margins from 1800 twips (90pt) to 1440 twips (72pt).
$ps = $d.PageSetup
$ps.LeftMargin = 72
$ps.RightMargin = 72

Write-Output "Done. Paragraphs=$($d.Paragraphs.Count)"
